# MN: Proyectos TestNG finales web and mobile
#
# - crearCuenta!E2 value changes from "No" to "Si"
# - Selections move on "calculadora" and "alertsTools"
# - The active/selected sheet moves from "calculadora" to "date",
#   with a new selection on "date" too.

$wb = $excel.ActiveWorkbook

# crearCuenta: E2 "No" -> "Si"
$wsCrearCuenta = $wb.Worksheets.Item("crearCuenta")
$wsCrearCuenta.Range("E2").Value = "Si"

# alertsTools: selection moves from E4 to D2 (sheet stays inactive)
$wsAlertsTools = $wb.Worksheets.Item("alertsTools")
$wsAlertsTools.Activate()
$wsAlertsTools.Range("D2").Select()

# calculadora: selection moves from H2 to C2 (sheet loses tabSelected)
$wsCalculadora = $wb.Worksheets.Item("calculadora")
$wsCalculadora.Activate()
$wsCalculadora.Range("C2").Select()

# date: becomes the active/selected sheet, selection moves from C11 to E16
$wsDate = $wb.Worksheets.Item("date")
$wsDate.Activate()
$wsDate.Range("E16").Select()
